$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in C9 (appearing -> appeari)
$ws.Range("C9").Value = "adding hidden keywords is a black-hat technique to improve SEO, what could be easily recognized by SEO robots, if so, website could be penalized and not appeari in indexing. Also it might be affecting accessibility, since the keywords might be getting focus or be read by screen readers"

# Fill new row 11 content
$ws.Range("A11").Value = "Both"
$ws.Range("B11").Value = "keywords in alt attribute in img elements"
$ws.Range("C11").Value = "1. Alt attribute is meant to containt a short description of the content of the image, so screen readers can read them out for accessibillity purposes"
$ws.Range("D11").Value = "Always use short description of the content of each image on the website"

# Fill new row 12 content
$ws.Range("C12").Value = "2. Placing key words in alt attribute is a black hat technique, if found by SEO robots, website could be penalized and not appear in indexing"

# Merge the cells that span rows 11-12
$ws.Range("A11:A12").Merge()
$ws.Range("B11:B12").Merge()
$ws.Range("D11:D12").Merge()

# Center alignment + wrap text for the new block (rows 11-12)
$block = $ws.Range("A11:D12")
$block.HorizontalAlignment = -4108
$block.VerticalAlignment = -4108
$block.WrapText = $true

# Row heights
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 54.75

# Update the active selection to match the author's last cursor position
$null = $ws.Range("D8").Select()
